$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bank account upload flags from "Y" to "N" for columns AR through AX on row 2
$ws.Range("AR2:AX2").Value = "N"

# Update the view: scroll so AR1 is top-left, and select AT19
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 44
$ws.Range("AT19").Select()
